$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'38.374.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "'2.081.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'227.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("D7").Value = "'60.58"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.380"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").Value = "'0.0834"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "'2.391.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").Value = "'14.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").Value = "'22.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.30%  "
$ws.Range("D15").Value = "'0.784"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "'5.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.75%  "
$ws.Range("D17").Value = "'2.070.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "'38.278.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").Value = "'71.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.23%  "
$ws.Range("D20").Value = "'6.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("D22").Value = "'225.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D24").Value = "'2.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").Value = "'169.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("D27").Value = "'9.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("E28").Value = "  +5.12%  "
$ws.Range("D29").Value = "'19.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("E30").Value = "  +8.39%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  +5.79%  "
$ws.Range("D33").Value = "'4.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.61%  "
$ws.Range("D34").Value = "'4.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("D35").Value = "'0.0605"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").Value = "'6.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("E37").Value = "  +1.65%  "
$ws.Range("E38").Value = "  +3.11%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("D41").Value = "'1.536.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").Value = "'99.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.34%  "
$ws.Range("D43").Value = "'0.0219"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").Value = "'0.0930"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").Value = "'2.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").Value = "'7.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.08%  "
$ws.Range("D47").Value = "'4.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").Value = "'1.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("D51").Value = "'2.284.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.12%  "
